$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '320.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.31%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.56%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.908'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '12.77%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08015'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.11%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.575'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.28%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.643'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.53%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.873'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.86%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.943'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.95%'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9329'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.41%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1242'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.70%'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1952'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.32%'
$ws.Range("B13").Value = 'MCDex'
$ws.Range("C13").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.762'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '29.87%'
$ws.Range("B14").Value = 'MandalaExchangeToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09154'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.32%'
$ws.Range("B15").Value = 'BitrueCoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03463'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.93%'
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.09597'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.58%'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001295'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-7.86%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006091'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-5.77%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.358'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.02%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.01%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1412'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '6.30%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2416'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.48%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04475'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.85%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001267'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.63%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004374'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.47%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001144'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-11.32%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.05%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02402'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.48%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05183'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.43%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007482'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.76%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1405'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.75%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.009121'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.01%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002099'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.53%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01151'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '41.05%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006761'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.79%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000753'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.32%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003017'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '5.79%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-42.66%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002107'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.32%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002007'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.32%'
